# Auto-generated PowerShell COM-interop script to apply the
# market-data refresh described by the commit "chore: update Sheets via scheduled runner".
# Each row updates cached currentAveragePrice*/LevePrice*/LeveProfit* columns (H-N);
# some rows drop the LeveProfitNQ (M) cell entirely when LevePriceNQ (K) becomes 0,
# matching the existing convention used throughout the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 150.54546
$ws.Range("I33").Value = 143.78947
$ws.Range("J33").Value = 193.33333
$ws.Range("K33").Value = 143.78947
$ws.Range("L33").Value = 193.33333
$ws.Range("M33").Value = 85.21053000000001
$ws.Range("N33").Value = -651.3333299999999

$ws.Range("H74").Value = 11114345
$ws.Range("I74").Value = 11114345
$ws.Range("K74").Value = 11114345
$ws.Range("M74").Value = -11113409

$ws.Range("H77").Value = 11114345
$ws.Range("I77").Value = 11114345
$ws.Range("K77").Value = 55571725
$ws.Range("M77").Value = -55567045

$ws.Range("H129").Value = 956.925
$ws.Range("J129").Value = 1134.9678
$ws.Range("L129").Value = 3404.9034
$ws.Range("N129").Value = -13404.9034

$ws.Range("H132").Value = 20082130
$ws.Range("I132").Value = 25101702
$ws.Range("J132").Value = 3844
$ws.Range("K132").Value = 75305106
$ws.Range("L132").Value = 11532
$ws.Range("M132").Value = -75302576
$ws.Range("N132").Value = -16592

$ws.Range("H137").Value = 3192
$ws.Range("I137").Value = 2576.2666
$ws.Range("J137").Value = 6655.5
$ws.Range("K137").Value = 7728.7998
$ws.Range("L137").Value = 19966.5
$ws.Range("M137").Value = -5178.7998
$ws.Range("N137").Value = -25066.5

$ws.Range("H138").Value = 2825.9062
$ws.Range("I138").Value = 2101.2632
$ws.Range("J138").Value = 3004.7144
$ws.Range("K138").Value = 6303.7896
$ws.Range("L138").Value = 9014.143199999999
$ws.Range("M138").Value = -1163.7896
$ws.Range("N138").Value = -19294.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3949.75
$ws.Range("I2").Value = 2250.5
$ws.Range("K2").Value = 2250.5
$ws.Range("M2").Value = -2137.5

$ws.Range("H32").Value = 11519.319
$ws.Range("I32").Value = 7880.477
$ws.Range("K32").Value = 7880.477
$ws.Range("M32").Value = -7593.477

$ws.Range("H110").Value = 958.2857
$ws.Range("I110").Value = 936
$ws.Range("J110").Value = 998.4
$ws.Range("K110").Value = 936
$ws.Range("L110").Value = 998.4
$ws.Range("M110").Value = 1109
$ws.Range("N110").Value = -5088.4

$ws.Range("H116").Value = 3949.75
$ws.Range("I116").Value = 2250.5
$ws.Range("K116").Value = 2250.5
$ws.Range("M116").Value = 43.5

$ws.Range("H132").Value = 5742
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 5742
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 17226
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -22286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3949.75
$ws.Range("I3").Value = 2250.5
$ws.Range("K3").Value = 2250.5
$ws.Range("M3").Value = -2136.5

$ws.Range("H86").Value = 1726.1538
$ws.Range("I86").Value = 1658.4445
$ws.Range("J86").Value = 1878.5
$ws.Range("K86").Value = 1658.4445
$ws.Range("L86").Value = 1878.5
$ws.Range("M86").Value = -535.4445000000001
$ws.Range("N86").Value = -4124.5

$ws.Range("H89").Value = 1726.1538
$ws.Range("I89").Value = 1658.4445
$ws.Range("J89").Value = 1878.5
$ws.Range("K89").Value = 8292.2225
$ws.Range("L89").Value = 9392.5
$ws.Range("M89").Value = -2676.2225
$ws.Range("N89").Value = -20624.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2998.75
$ws.Range("I31").Value = 1326.6
$ws.Range("J31").Value = 5785.6665
$ws.Range("K31").Value = 1326.6
$ws.Range("L31").Value = 5785.6665
$ws.Range("M31").Value = -1031.6
$ws.Range("N31").Value = -6375.6665

$ws.Range("H34").Value = 2998.75
$ws.Range("I34").Value = 1326.6
$ws.Range("J34").Value = 5785.6665
$ws.Range("K34").Value = 1326.6
$ws.Range("L34").Value = 5785.6665
$ws.Range("M34").Value = -1124.6
$ws.Range("N34").Value = -6189.6665

$ws.Range("H99").Value = 9250
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 9250
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 9250
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -12246

$ws.Range("H126").Value = 9250
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 9250
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 27750
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -32690

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5544.684
$ws.Range("I3").Value = 3997.182
$ws.Range("J3").Value = 7672.5
$ws.Range("K3").Value = 11991.546
$ws.Range("L3").Value = 23017.5
$ws.Range("M3").Value = -11879.546
$ws.Range("N3").Value = -23241.5

$ws.Range("H60").Value = 10611.182
$ws.Range("I60").Value = 125.666664
$ws.Range("J60").Value = 14543.25
$ws.Range("K60").Value = 376.999992
$ws.Range("L60").Value = 43629.75
$ws.Range("M60").Value = -125.999992
$ws.Range("N60").Value = -44131.75

$ws.Range("H92").Value = 62504184
$ws.Range("I92").Value = 553.3333
$ws.Range("J92").Value = 100006360
$ws.Range("K92").Value = 1659.9999
$ws.Range("L92").Value = 300019080
$ws.Range("M92").Value = -411.9999
$ws.Range("N92").Value = -300021576

$ws.Range("H113").Value = 763.5263
$ws.Range("I113").Value = 685.3043
$ws.Range("J113").Value = 883.4666999999999
$ws.Range("K113").Value = 2055.9129
$ws.Range("L113").Value = 2650.4001
$ws.Range("M113").Value = 114.0870999999997
$ws.Range("N113").Value = -6990.4001

$ws.Range("H131").Value = 11628859
$ws.Range("J131").Value = 1055.1111
$ws.Range("L131").Value = 3165.3333
$ws.Range("N131").Value = -13245.3333

$ws.Range("H137").Value = 3453.3333
$ws.Range("I137").Value = 3184
$ws.Range("J137").Value = 4800
$ws.Range("K137").Value = 9552
$ws.Range("L137").Value = 14400
$ws.Range("M137").Value = -4452
$ws.Range("N137").Value = -24600

$ws.Range("H138").Value = 2355.4
$ws.Range("I138").Value = 1407.8572
$ws.Range("K138").Value = 4223.571599999999
$ws.Range("M138").Value = 916.4284000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3212.9473
$ws.Range("I122").Value = 2549.2144
$ws.Range("K122").Value = 7647.6432
$ws.Range("M122").Value = -5197.6432

$ws.Range("H126").Value = 3943.8403
$ws.Range("I126").Value = 2902.1133
$ws.Range("J126").Value = 5290.4634
$ws.Range("K126").Value = 8706.339899999999
$ws.Range("L126").Value = 15871.3902
$ws.Range("M126").Value = -6236.339899999999
$ws.Range("N126").Value = -20811.3902

$ws.Range("H132").Value = 3012.739
$ws.Range("I132").Value = 1277.0834
$ws.Range("K132").Value = 3831.2502
$ws.Range("M132").Value = -1301.2502

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3657.182
$ws.Range("I136").Value = 1930.1875
$ws.Range("J136").Value = 5282.5884
$ws.Range("K136").Value = 5790.5625
$ws.Range("L136").Value = 15847.7652
$ws.Range("M136").Value = -3240.5625
$ws.Range("N136").Value = -20947.7652
